$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row at position 2 for the updated "A 34417-2023" entry ---
# (its previous, out-of-date entry further down the sheet will be removed below)
$ws.Rows("2:2").Insert()

# --- 2) Populate the new row 2 with the refreshed data for "A 34417-2023" ---
$ws.Range("A2").Value = "A 34417-2023"

$ws.Range("B2").Value = 45139
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"

$ws.Range("C2").Value = 45177
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"

$ws.Range("D2").Value = "STOCKHOLMS LÄN"
$ws.Range("E2").Value = "HANINGE"

$ws.Range("G2").Value = 3.8
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 20

$ws.Range("R2").Value = "Porslinsblå spindling`r`nBarrviolspindling`r`nDvärgbägarlav`r`nLeptoporus erubescens`r`nSpillkråka`r`nVedtrappmossa`r`nBronshjon`r`nDropptaggsvamp`r`nFällmossa`r`nGrön sköldmossa`r`nGuldlockmossa`r`nKornknutmossa`r`nRödgul trumpetsvamp`r`nStubbspretmossa`r`nSvavelriska`r`nSårläka`r`nVågbandad barkbock`r`nVanlig groda`r`nBlåsippa`r`nRevlummer"
$ws.Range("R2").WrapText = $true

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/artfynd/A 34417-2023.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/kartor/A 34417-2023.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomål/A 34417-2023.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomålsmail/A 34417-2023.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsyn/A 34417-2023.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsynsmail/A 34417-2023.docx")'

# Row height is fixed at 15 throughout the sheet; setting it after the wrapped,
# multi-line R2 text avoids Excel's auto-fit bumping it up.
$ws.Rows("2:2").RowHeight = 15

# --- 3) Remove the old (now stale) "A 34417-2023" row further down the sheet ---
# After the insert above it has shifted from row 11 to row 12.
$ws.Rows("12:12").Delete()

# --- 4) Every remaining record's "Förändrad" date moved from 2023-09-06 to 2023-09-08 ---
$ws.Range("C2:C91").Value = 45177
$ws.Range("C2:C91").NumberFormat = "YYYY-MM-DD"

Write-Output "done"
